# Update cached/computed values in columns C:F for rows 3-10 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = @(32892.48046875001, 2183.59375, 1441.171875, 694.5902301087816)
    4  = @(58135.54687500001, 3859.375, 2547.1875, 1227.647848564358)
    5  = @(75729.19921875001, 5027.34375, 3318.046875, 1599.17285536673)
    6  = @(85673.43750000001, 5687.5, 3753.75, 1809.165250515896)
    7  = @(87968.26171875001, 5839.84375, 3854.296875, 1857.625034011858)
    8  = @(82613.67187500001, 5484.375, 3619.6875, 1744.552205854614)
    9  = @(69609.66796875001, 4621.09375, 3049.921875, 1469.946766044166)
    10 = @(48956.25000000001, 3250, 2145, 1033.808714580512)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Range("C$row").Value = $rowValues[0]
    $ws.Range("D$row").Value = $rowValues[1]
    $ws.Range("E$row").Value = $rowValues[2]
    $ws.Range("F$row").Value = $rowValues[3]
}
